$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the formatting used by the
# other header cells (copy G1's format, e.g. bold font/border/alignment).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill the new column's data rows with 0 (unstyled, like the other data cells).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
